$wb = $excel.ActiveWorkbook

# ===== Sheet: Overview =====
$ws = $wb.Worksheets.Item("Overview")

# -- update cell values --
$ws.Range("A2").Value = "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-22 06:53:26"

$ws.Range("A3").Value = "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-22 06:53:26"

$ws.Range("A4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-22 06:56:57"

# -- rebuild hyperlinks (delete all, then re-add in order so rIds line up) --
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/5f868707-b1fc-49bf-888a-1db5ffc5e40b.md", "", "", "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md", "", "", "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md") | Out-Null

# ===== Sheet: zh-cn =====
$ws = $wb.Worksheets.Item("zh-cn")

# -- update cell values --
$ws.Range("A2").Value = "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-22 06:53:18"
$ws.Range("F2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.md"
$ws.Range("G2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-22 06:53:59"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-22 06:53:18"
$ws.Range("F3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.md"
$ws.Range("G3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-22 06:53:59"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-22 06:56:49"
$ws.Range("F4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md"
$ws.Range("G4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-22 06:55:56"
$ws.Range("J4").Value = "Include"

# -- rebuild hyperlinks (delete all, then re-add in order so rIds line up) --
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/5f868707-b1fc-49bf-888a-1db5ffc5e40b.md", "", "", "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78ddd65bb8c49b6bcf24596255f5eeb937bb05dd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/e875070a273520fbfe12c80eb0ed65a06c0d671a/e2e/5f868707-b1fc-49bf-888a-1db5ffc5e40b.md", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c88603e37eeb464f8b4b02ccf19912ad7252e74b/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md", "", "", "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25bb827ec779a78c5325a01e8ae547469650ce62/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/2c55268b1c19e7e198ec6b6d58d619186abd2397/e2e/1761058d-58d9-4c45-a20e-e70b262a33d9.md", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/adf8019aca1fc5bf9d296aff5e3f5e9aeacb9001/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25bb827ec779a78c5325a01e8ae547469650ce62/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/2c55268b1c19e7e198ec6b6d58d619186abd2397/e2e/1761058d-58d9-4c45-a20e-e70b262a33d9.md", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/adf8019aca1fc5bf9d296aff5e3f5e9aeacb9001/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.zh-cn.xlf") | Out-Null

# ===== Sheet: de-de =====
$ws = $wb.Worksheets.Item("de-de")

# -- update cell values --
$ws.Range("A2").Value = "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"
$ws.Range("E2").Value = "2016-03-22 06:53:26"
$ws.Range("F2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.md"
$ws.Range("G2").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"
$ws.Range("H2").Value = "2016-03-22 06:54:13"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"
$ws.Range("E3").Value = "2016-03-22 06:53:26"
$ws.Range("F3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.md"
$ws.Range("G3").Value = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"
$ws.Range("H3").Value = "2016-03-22 06:54:13"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf"
$ws.Range("E4").Value = "2016-03-22 06:56:57"
$ws.Range("F4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md"
$ws.Range("G4").Value = "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf"
$ws.Range("H4").Value = "2016-03-22 06:56:11"
$ws.Range("J4").Value = "Include"

# -- rebuild hyperlinks (delete all, then re-add in order so rIds line up) --
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/5f868707-b1fc-49bf-888a-1db5ffc5e40b.md", "", "", "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3f3d2621bab3153d1ffd44c17c1899530266fbb/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/f3390c930c3418c7436433d298026df4a0ba65ba/e2e/5f868707-b1fc-49bf-888a-1db5ffc5e40b.md", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1be1809e585117ce4c9638aa03664843888d8bfb/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c.md", "", "", "ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56abd67a8b8fdabb4e80ed514086018986e412e2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/361f9022d354657eca286d6ae9983daac72d8146/e2e/1761058d-58d9-4c45-a20e-e70b262a33d9.md", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc730f15ea71c6621c3e88fff2c8651b83603bbe/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf", "", "", "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e8835e41e57cc90edf540aac9441a5cc33320292/e2e/ffffff2f4f3802-0c9a-43e9-8aa4-1743f50736c5.md", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56abd67a8b8fdabb4e80ed514086018986e412e2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/361f9022d354657eca286d6ae9983daac72d8146/e2e/1761058d-58d9-4c45-a20e-e70b262a33d9.md", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc730f15ea71c6621c3e88fff2c8651b83603bbe/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf", "", "", "5f868707-b1fc-49bf-888a-1db5ffc5e40b.824ba228a2a09edf8631909b229db69fac306d0a.de-de.xlf") | Out-Null

